$d = $word.ActiveDocument

# The final paragraph of the document is currently:
#   "Le fait de remplacer le StringBuffer par un BigInteger a ete determinant."
# It gets superseded below by a longer block of new content that ends by
# re-stating that very same sentence.
$parLast = $d.Paragraphs.Last

# Remove the "_GoBack" bookmark from its current (old) location -- the empty
# paragraph just before $parLast. It will be re-created at the end of the
# new content block being inserted below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Build the block of new paragraphs (the expanded write-up plus the
# relocated "Le fait de remplacer..." paragraph) as literal WordprocessingML.
$newParagraphsXml = @"
<w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Le fait de remplacer le StringBuffer par un BigInteger a été déterminant.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>On inclut le module dans le projet principal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>On s’occupe de la déconnexion</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:i/><w:iCs/><w:color w:val="808080"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Problème : j’utilisais </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="808080"/></w:rPr><w:t>informUser = R.string.Password_failure;</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> pour récupérer une chaine dans les ressources au sein du code java, ce qui déclenchait une erreur car </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="808080"/></w:rPr><w:t>R.string.Password_failure</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="808080"/></w:rPr><w:t xml:space="preserve"> est considéré comme un entier par java. J’ai donc plutôt utilisé </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFE4FF"/></w:rPr><w:t>informUser</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> =</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>getResources().getString(R.string.</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="660E7A"/></w:rPr><w:t>Password_failure</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PrformatHTML"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>ce qui fait bien ce que je désirais, c’est-à-dire mettre le string « Password_failure » dans informUser.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

# Insert the new paragraphs in place of the final paragraph: collapsing the
# range to its start and inserting WordprocessingML there replaces that
# paragraph's own content with the first inserted paragraph, while the
# remaining inserted paragraphs are appended as brand-new paragraphs after
# it -- which is exactly the "Le fait de remplacer..." paragraph being
# replaced by the full new block (whose own last paragraph re-creates that
# same sentence, now carrying the relocated bookmark).
$insertionRange = $parLast.Range
$insertionRange.Collapse(1)
$insertionRange.InsertXML($newParagraphsXml)
